$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "30.664.32"
$ws.Range("E2").Value = "  +1.22%  "

$ws.Range("D3").Value = "1.871.03"
$ws.Range("E3").Value = "  +0.55%  "

Set-TextValue $ws.Range("D4") "0.9998"
$ws.Range("E4").Value = "  -0.09%  "

Set-TextValue $ws.Range("D5") "235.85"
$ws.Range("E5").Value = "  +1.11%  "

Set-TextValue $ws.Range("D6") "0.9992"
$ws.Range("E6").Value = "  -0.14%  "

Set-TextValue $ws.Range("D7") "0.4717"
$ws.Range("E7").Value = "  -0.49%  "

Set-TextValue $ws.Range("D8") "0.2776"
$ws.Range("E8").Value = "  +1.33%  "

Set-TextValue $ws.Range("D9") "0.06401"
$ws.Range("E9").Value = "  -0.40%  "

Set-TextValue $ws.Range("D10") "18.18"
$ws.Range("E10").Value = "  +12.32%  "

$ws.Range("D11").Value = "1.870.89"
$ws.Range("E11").Value = "  +0.51%  "

Set-TextValue $ws.Range("D12") "0.07444"
$ws.Range("E12").Value = "  +0.32%  "

Set-TextValue $ws.Range("D13") "4.992"
$ws.Range("E13").Value = "  -0.13%  "

Set-TextValue $ws.Range("D14") "85.60"
$ws.Range("E14").Value = "  +0.07%  "

Set-TextValue $ws.Range("D15") "0.6416"
$ws.Range("E15").Value = "  +1.79%  "

$ws.Range("D16").Value = "30.639.44"
$ws.Range("E16").Value = "  +1.32%  "

Set-TextValue $ws.Range("D17") "246.78"
$ws.Range("E17").Value = "  +6.62%  "

Set-TextValue $ws.Range("D19") "12.92"
$ws.Range("E19").Value = "  +0.97%  "

Set-TextValue $ws.Range("D20") "0.000007432"
$ws.Range("E20").Value = "  +1.13%  "

Set-TextValue $ws.Range("D21") "0.9991"
$ws.Range("E21").Value = "  -0.25%  "

Set-TextValue $ws.Range("D22") "5.005"
$ws.Range("E22").Value = "  -1.40%  "

Set-TextValue $ws.Range("D23") "6.146"
$ws.Range("E23").Value = "  +2.70%  "

Set-TextValue $ws.Range("D24") "9.404"
$ws.Range("E24").Value = "  +1.49%  "

Set-TextValue $ws.Range("D25") "164.91"
$ws.Range("E25").Value = "  -1.20%  "

Set-TextValue $ws.Range("D26") "18.49"
$ws.Range("E26").Value = "  +3.90%  "

Set-TextValue $ws.Range("D27") "1.902"
$ws.Range("E27").Value = "  +1.63%  "

Set-TextValue $ws.Range("D28") "0.1022"
$ws.Range("E28").Value = "  +2.55%  "

Set-TextValue $ws.Range("D30") "4.110"
$ws.Range("E30").Value = "  -1.65%  "

Set-TextValue $ws.Range("D31") "3.885"
$ws.Range("E31").Value = "  -0.91%  "

Set-TextValue $ws.Range("D32") "0.04941"
$ws.Range("E32").Value = "  +0.92%  "

Set-TextValue $ws.Range("D33") "1.165"
$ws.Range("E33").Value = "  +1.92%  "

Set-TextValue $ws.Range("D34") "0.7173"
$ws.Range("E34").Value = "  +0.55%  "

Set-TextValue $ws.Range("D35") "2.704"
$ws.Range("E35").Value = "  +0.17%  "

Set-TextValue $ws.Range("D36") "0.01910"
$ws.Range("E36").Value = "  -0.09%  "

Set-TextValue $ws.Range("D37") "2.694"
$ws.Range("E37").Value = "  +2.09%  "

Set-TextValue $ws.Range("D38") "0.8823"
$ws.Range("E38").Value = "  -2.18%  "

Set-TextValue $ws.Range("D39") "2.009"
$ws.Range("E39").Value = "  +1.65%  "

Set-TextValue $ws.Range("D40") "106.08"
$ws.Range("E40").Value = "  +0.38%  "

Set-TextValue $ws.Range("D41") "0.9988"
$ws.Range("E41").Value = "  -0.17%  "

Set-TextValue $ws.Range("D42") "0.4127"
$ws.Range("E42").Value = "  +0.64%  "

Set-TextValue $ws.Range("D43") "5.583"
$ws.Range("E43").Value = "  +0.39%  "

Set-TextValue $ws.Range("D44") "7.419"
$ws.Range("E44").Value = "  +5.33%  "

Set-TextValue $ws.Range("D45") "62.45"
$ws.Range("E45").Value = "  +2.22%  "

Set-TextValue $ws.Range("D47") "8.740"
$ws.Range("E47").Value = "  -0.26%  "

Set-TextValue $ws.Range("D48") "33.80"
$ws.Range("E48").Value = "  +2.44%  "

Set-TextValue $ws.Range("D51") "0.3725"
$ws.Range("E51").Value = "  +0.87%  "

$ws.Range("E18").Value = "  -0.03%  "

$ws.Range("E29").Value = "  -0.37%  "

$ws.Range("E46").Value = "  +2.75%  "

# Row 49: Cronos -> NEARProtocol
$ws.Range("B49").Value = "NEARProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D49") "1.389"
$ws.Range("E49").Value = "  -0.83%  "

# Row 50: NEARProtocol -> Cronos
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws.Range("D50") "0.05566"
$ws.Range("E50").Value = "  -0.26%  "
